# Fixed a bug in randomMoveSymbols
#
# The data rows (2-25, columns A:F) got shuffled into a new row order.
# Capture the current A:F block for each data row, then write it back
# out in the permuted order below (destination row -> source row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcForDest = @{
    2  = 12
    3  = 10
    4  = 15
    5  = 9
    6  = 3
    7  = 4
    8  = 11
    9  = 6
    10 = 13
    11 = 14
    12 = 5
    13 = 8
    14 = 2
    15 = 7
    16 = 20
    17 = 21
    18 = 18
    19 = 17
    20 = 16
    21 = 19
    22 = 23
    23 = 22
    24 = 24
    25 = 25
}

# Snapshot the original values of rows 2-25, columns A:F, before overwriting.
$original = @{}
for ($r = 2; $r -le 25; $r++) {
    $original[$r] = $ws.Range("A$r`:F$r").Value2
}

# Write each destination row using the snapshot of its mapped source row.
for ($r = 2; $r -le 25; $r++) {
    $srcRow = $srcForDest[$r]
    $ws.Range("A$r`:F$r").Value2 = $original[$srcRow]
}
